$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark from the "Entrega" date cell
#    ("24/07/16") in the "Ciclo de Vida" table. Bookmarks.Item(...).Delete()
#    does not actually drop the special "_GoBack" bookmark in this host, so
#    we rewrite the two characters spanning it ("4/") which forces the
#    bookmark to be dropped when the paragraph's runs are rebuilt.
# ------------------------------------------------------------------
$goBackOld = $d.Content.Duplicate
$goBackOld.Find.Execute("24/07/16")
if ($goBackOld.Find.Found) {
    $span = $d.Range($goBackOld.Start + 1, $goBackOld.Start + 3)
    $spanText = $span.Text
    $span.Delete()
    $reins = $d.Range($goBackOld.Start + 1, $goBackOld.Start + 1)
    $reins.InsertBefore($spanText)
}

# ------------------------------------------------------------------
# 2) Add the new row to the "Ciclo de Vida da Ordem de Serviço" table
#    (table #5) describing the 19/08/16 publication event.
# ------------------------------------------------------------------
$t = $d.Tables.Item(5)
$newRow = $t.Rows.Add()

# -- Cell 1: Date, centered --------------------------------------------------
$cell1 = $newRow.Cells.Item(1)
$cell1.Range.Text = "19/08/16"
$cell1.Range.ParagraphFormat.Alignment = 1

# -- Cell 2: Event description (3 paragraphs) --------------------------------
$cell2 = $newRow.Cells.Item(2)
$cell2.Range.Text = "Publicação de novas versões dos casos de usos:" + "`r" + `
    "ARRUC0910 - Parametrizar Bancos e" + "`r" + `
    "ARRUC0940 - Parametrizar Plano de Contas, e também do dicionário de dados e modelo de casos de usos."

# Split "Publicação de novas versões d" | "os casos de usos:" into two runs
$r1 = $d.Content.Duplicate
$r1.Find.Execute("os casos de usos:")
if ($r1.Find.Found) {
    $r1.Font.Bold = 1
    $r1.Font.Bold = 0
}

# Split "ARRUC0910 - Parametrizar Bancos" | " e" into two runs
$r2 = $d.Content.Duplicate
$r2.Find.Execute("ARRUC0910 - Parametrizar Bancos")
if ($r2.Find.Found) {
    $r2.Collapse(0)
    $r2.MoveEnd(1, 2)
    $r2.Font.Bold = 1
    $r2.Font.Bold = 0
}

# Re-create the "_GoBack" bookmark right after "ARRUC0940 - Parametrizar
# Plano de Contas" (this also naturally splits that run from the trailing
# ", e também ..." text).
$r3 = $d.Content.Duplicate
$r3.Find.Execute("ARRUC0940 - Parametrizar Plano de Contas")
if ($r3.Find.Found) {
    $r3.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $r3)
}

# -- Cell 3: Executor ----------------------------------------------------
$cell3 = $newRow.Cells.Item(3)
$cell3.Range.Text = "CIAT - João Paulo"
